{"js": "// Remove the comma after \"curious\" in the \"What are some limitations of\n// this dataset?\" answer paragraph:\n//   \"It is very curious, why some of the projects ...\"\n// becomes\n//   \"It is very curious why some of the projects ...\"\nconst body = context.document.body;\n\nconst searchResults = body.search(\"curious, why\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Expected text \"curious, why\" not found in document body.');\n}\n\n// Replace the matched text in place; Word keeps the surrounding run\n// formatting (the 2B2B2B font color) intact automatically.\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"curious why\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Remove the comma after \"curious\" in the \"What are some limitations of\n# this dataset?\" answer paragraph:\n#   \"It is very curious, why some of the projects ...\"\n# becomes\n#   \"It is very curious why some of the projects ...\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"curious, why\",  # FindText\n    $false,          # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    1,               # Wrap (wdFindContinue)\n    $false,          # Format\n    \"curious why\",   # ReplaceWith\n    2                # Replace (wdReplaceAll)\n)\n\nif (-not $found) {\n    throw 'Expected text \"curious, why\" not found in document content.'\n}\n"}
